# Auto-generated script applying the Kujata_Profits market-data refresh
# described in the commit "chore: update Sheets via scheduled runner".
$wb = $excel.ActiveWorkbook

# --- Sheet ALC, hunk 0 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value2 = 1001.9545
$ws.Range("J17").Value2 = 1001.9545
$ws.Range("L17").Value2 = 3005.8635
$ws.Range("N17").Value2 = -3341.8635

# --- Sheet ALC, hunk 1 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value2 = 367.625
$ws.Range("I96").Value2 = 291.57144
$ws.Range("K96").Value2 = 874.71432
$ws.Range("M96").Value2 = 498.28568

# --- Sheet ALC, hunk 2 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value2 = 1946.0
$ws.Range("I100").Value2 = 1736.0
$ws.Range("J100").Value2 = 2996.0
$ws.Range("K100").Value2 = 1736.0
$ws.Range("L100").Value2 = 2996.0
$ws.Range("M100").Value2 = -1195.0
$ws.Range("N100").Value2 = -4078.0

# --- Sheet ALC, hunk 3 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value2 = 613.0455
$ws.Range("I103").Value2 = 492.125
$ws.Range("J103").Value2 = 682.1429
$ws.Range("K103").Value2 = 1476.375
$ws.Range("L103").Value2 = 2046.4287
$ws.Range("M103").Value2 = -890.375
$ws.Range("N103").Value2 = -3218.4287

# --- Sheet ALC, hunk 4 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value2 = 1933.65
$ws.Range("J112").Value2 = 1933.65
$ws.Range("L112").Value2 = 5800.950000000001
$ws.Range("N112").Value2 = -8016.950000000001

# --- Sheet ALC, hunk 5 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value2 = 1911.7949
$ws.Range("I125").Value2 = 1674.909
$ws.Range("J125").Value2 = 2218.353
$ws.Range("K125").Value2 = 15074.181
$ws.Range("L125").Value2 = 19965.177
$ws.Range("M125").Value2 = -12614.181
$ws.Range("N125").Value2 = -24885.177

# --- Sheet ALC, hunk 6 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value2 = 28574154.0
$ws.Range("I135").Value2 = 966.84
$ws.Range("J135").Value2 = 100007120.0
$ws.Range("K135").Value2 = 8701.56
$ws.Range("L135").Value2 = 900064080.0
$ws.Range("M135").Value2 = -6166.559999999999
$ws.Range("N135").Value2 = -900069150.0

# --- Sheet ALC, hunk 7 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value2 = 1090.6129
$ws.Range("I137").Value2 = 740.5909
$ws.Range("J137").Value2 = 1946.2222
$ws.Range("K137").Value2 = 2221.7727
$ws.Range("L137").Value2 = 5838.6666
$ws.Range("M137").Value2 = 328.2273
$ws.Range("N137").Value2 = -10938.6666

# --- Sheet ALC, hunk 8 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value2 = 2941.817
$ws.Range("J138").Value2 = 2942.2952
$ws.Range("L138").Value2 = 8826.8856
$ws.Range("N138").Value2 = -19106.8856

# --- Sheet ARM, hunk 9 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 14686.172
$ws.Range("I32").Value2 = 9830.862
$ws.Range("J32").Value2 = 22732.115
$ws.Range("K32").Value2 = 9830.862
$ws.Range("L32").Value2 = 22732.115
$ws.Range("M32").Value2 = -9543.862
$ws.Range("N32").Value2 = -23306.115

# --- Sheet ARM, hunk 10 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value2 = 29800.0
$ws.Range("J44").Value2 = 29800.0
$ws.Range("L44").Value2 = 29800.0
$ws.Range("N44").Value2 = -30776.0

# --- Sheet ARM, hunk 11 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value2 = 32409.6

# --- Sheet ARM, hunk 12 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value2 = 1471.825
$ws.Range("I74").Value2 = 743.86957
$ws.Range("J74").Value2 = 2456.7058
$ws.Range("K74").Value2 = 743.86957
$ws.Range("L74").Value2 = 2456.7058
$ws.Range("M74").Value2 = 130.13043
$ws.Range("N74").Value2 = -4204.7058

# --- Sheet ARM, hunk 13 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value2 = 1471.825
$ws.Range("I77").Value2 = 743.86957
$ws.Range("J77").Value2 = 2456.7058
$ws.Range("K77").Value2 = 3719.34785
$ws.Range("L77").Value2 = 12283.529
$ws.Range("M77").Value2 = 648.6521500000003
$ws.Range("N77").Value2 = -21019.529

# --- Sheet ARM, hunk 14 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value2 = 6782.75
$ws.Range("I97").Value2 = 584.2727
$ws.Range("K97").Value2 = 584.2727
$ws.Range("M97").Value2 = -88.27269999999999

# --- Sheet ARM, hunk 15 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value2 = 8336138.0
$ws.Range("I102").Value2 = 8336138.0
$ws.Range("K102").Value2 = 8336138.0
$ws.Range("M102").Value2 = -8334516.0

# --- Sheet ARM, hunk 16 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value2 = 3783.0
$ws.Range("I122").Value2 = 3017.375
$ws.Range("J122").Value2 = 5008.0
$ws.Range("K122").Value2 = 9052.125
$ws.Range("L122").Value2 = 15024.0
$ws.Range("M122").Value2 = -6602.125
$ws.Range("N122").Value2 = -19924.0

# --- Sheet ARM, hunk 17 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value2 = 2286.0227
$ws.Range("I132").Value2 = 1852.4166
$ws.Range("J132").Value2 = 4237.25
$ws.Range("K132").Value2 = 5557.2498
$ws.Range("L132").Value2 = 12711.75
$ws.Range("M132").Value2 = -3027.2498
$ws.Range("N132").Value2 = -17771.75

# --- Sheet BSM, hunk 18 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value2 = 40375.57
$ws.Range("I97").Value2 = 9577.4
$ws.Range("K97").Value2 = 9577.4
$ws.Range("M97").Value2 = -8586.4

# --- Sheet BSM, hunk 19 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value2 = 50001196.0
$ws.Range("J99").Value2 = 1712.5
$ws.Range("L99").Value2 = 1712.5
$ws.Range("N99").Value2 = -4708.5

# --- Sheet BSM, hunk 20 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value2 = 30390.0
$ws.Range("J122").Value2 = 30390.0
$ws.Range("L122").Value2 = 30390.0
$ws.Range("N122").Value2 = -40190.0

# --- Sheet BSM, hunk 21 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value2 = 60000.0
$ws.Range("J126").Value2 = 60000.0
$ws.Range("L126").Value2 = 60000.0
$ws.Range("N126").Value2 = -69880.0

# --- Sheet BSM, hunk 22 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value2 = 4813.9644
$ws.Range("I134").Value2 = 1012.46155
$ws.Range("J134").Value2 = 8108.6
$ws.Range("K134").Value2 = 3037.38465
$ws.Range("L134").Value2 = 24325.8
$ws.Range("M134").Value2 = -502.38465
$ws.Range("N134").Value2 = -29395.8

# --- Sheet CRP, hunk 23 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 1695.6154
$ws.Range("I31").Value2 = 1663.44
$ws.Range("K31").Value2 = 1663.44
$ws.Range("M31").Value2 = -1368.44

# --- Sheet CRP, hunk 24 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value2 = 1695.6154
$ws.Range("I34").Value2 = 1663.44
$ws.Range("K34").Value2 = 1663.44
$ws.Range("M34").Value2 = -1461.44

# --- Sheet CRP, hunk 25 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value2 = 3055980.2
$ws.Range("I86").Value2 = 5575738.5
$ws.Range("J86").Value2 = 32270.8
$ws.Range("K86").Value2 = 5575738.5
$ws.Range("L86").Value2 = 32270.8
$ws.Range("M86").Value2 = -5574615.5
$ws.Range("N86").Value2 = -34516.8

# --- Sheet CRP, hunk 26 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value2 = 3055980.2
$ws.Range("I89").Value2 = 5575738.5
$ws.Range("J89").Value2 = 32270.8
$ws.Range("K89").Value2 = 27878692.5
$ws.Range("L89").Value2 = 161354.0
$ws.Range("M89").Value2 = -27873076.5
$ws.Range("N89").Value2 = -172586.0

# --- Sheet CRP, hunk 27 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value2 = 1237.0625
$ws.Range("I107").Value2 = 685.3
$ws.Range("K107").Value2 = 685.3
$ws.Range("M107").Value2 = 1234.7

# --- Sheet CRP, hunk 28 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value2 = 323547.3
$ws.Range("J141").Value2 = 323547.3
$ws.Range("L141").Value2 = 323547.3
$ws.Range("N141").Value2 = -333907.3

# --- Sheet CUL, hunk 29 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value2 = 174.72223
$ws.Range("I12").Value2 = 231.8
$ws.Range("J12").Value2 = 152.76923
$ws.Range("K12").Value2 = 695.4000000000001
$ws.Range("L12").Value2 = 458.30769
$ws.Range("M12").Value2 = -522.4000000000001
$ws.Range("N12").Value2 = -804.30769

# --- Sheet CUL, hunk 30 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value2 = 950.0
$ws.Range("J54").Value2 = 950.0
$ws.Range("L54").Value2 = 2850.0
$ws.Range("N54").Value2 = -3968.0

# --- Sheet CUL, hunk 31 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value2 = 5940.35
$ws.Range("J107").Value2 = 11516.3
$ws.Range("L107").Value2 = 34548.89999999999
$ws.Range("N107").Value2 = -38388.89999999999

# --- Sheet CUL, hunk 32 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value2 = 704.30554
$ws.Range("I113").Value2 = 590.0
$ws.Range("J113").Value2 = 736.9643
$ws.Range("K113").Value2 = 1770.0
$ws.Range("L113").Value2 = 2210.8929
$ws.Range("M113").Value2 = 400.0
$ws.Range("N113").Value2 = -6550.8929

# --- Sheet CUL, hunk 33 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value2 = 1023.86
$ws.Range("J122").Value2 = 1080.641
$ws.Range("L122").Value2 = 9725.769
$ws.Range("N122").Value2 = -14625.769

# --- Sheet CUL, hunk 34 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value2 = 23846132.0
$ws.Range("J131").Value2 = 43846.6
$ws.Range("L131").Value2 = 131539.8
$ws.Range("N131").Value2 = -141619.8

# --- Sheet CUL, hunk 35 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value2 = 27784416.0
$ws.Range("I137").Value2 = 46877100.0
$ws.Range("J137").Value2 = 13242.272
$ws.Range("K137").Value2 = 140631300.0
$ws.Range("L137").Value2 = 39726.81600000001
$ws.Range("M137").Value2 = -140626200.0
$ws.Range("N137").Value2 = -49926.81600000001

# --- Sheet CUL, hunk 36 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value2 = 25360.75
$ws.Range("J140").Value2 = 2851.4583
$ws.Range("L140").Value2 = 8554.374899999999
$ws.Range("N140").Value2 = -18914.3749

# --- Sheet GSM, hunk 37 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value2 = 0.0
$ws.Range("I17").Value2 = 0.0
$ws.Range("K17").Value2 = 0.0
$ws.Range("M17").ClearContents()

# --- Sheet GSM, hunk 38 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 17311384.0
$ws.Range("I70").Value2 = 14709623.0
$ws.Range("K70").Value2 = 14709623.0
$ws.Range("M70").Value2 = -14709353.0

# --- Sheet GSM, hunk 39 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value2 = 17311384.0
$ws.Range("I73").Value2 = 14709623.0
$ws.Range("K73").Value2 = 14709623.0
$ws.Range("M73").Value2 = -14708687.0

# --- Sheet GSM, hunk 40 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value2 = 1635.6666
$ws.Range("I122").Value2 = 1635.6666
$ws.Range("K122").Value2 = 4906.9998
$ws.Range("M122").Value2 = -2456.9998

# --- Sheet GSM, hunk 41 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value2 = 28935.166
$ws.Range("J134").Value2 = 28935.166
$ws.Range("L134").Value2 = 86805.498
$ws.Range("N134").Value2 = -91875.498

# --- Sheet LTW, hunk 42 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 1592.6923
$ws.Range("I22").Value2 = 1824.75
$ws.Range("J22").Value2 = 1489.5555
$ws.Range("K22").Value2 = 1824.75
$ws.Range("L22").Value2 = 1489.5555
$ws.Range("M22").Value2 = -1529.75
$ws.Range("N22").Value2 = -2079.5555

# --- Sheet LTW, hunk 43 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value2 = 1592.6923
$ws.Range("I27").Value2 = 1824.75
$ws.Range("J27").Value2 = 1489.5555
$ws.Range("K27").Value2 = 1824.75
$ws.Range("L27").Value2 = 1489.5555
$ws.Range("M27").Value2 = -1717.75
$ws.Range("N27").Value2 = -1703.5555

# --- Sheet LTW, hunk 44 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value2 = 399.9091
$ws.Range("I55").Value2 = 249.57143
$ws.Range("J55").Value2 = 663.0
$ws.Range("K55").Value2 = 249.57143
$ws.Range("L55").Value2 = 663.0
$ws.Range("M55").Value2 = -76.57142999999999
$ws.Range("N55").Value2 = -1009.0

# --- Sheet LTW, hunk 45 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value2 = 1032.3846
$ws.Range("I93").Value2 = 987.2
$ws.Range("J93").Value2 = 1183.0
$ws.Range("K93").Value2 = 987.2
$ws.Range("L93").Value2 = 1183.0
$ws.Range("M93").Value2 = 260.8
$ws.Range("N93").Value2 = -3679.0

# --- Sheet LTW, hunk 46 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value2 = 1465.8334
$ws.Range("I100").Value2 = 1198.75
$ws.Range("K100").Value2 = 1198.75
$ws.Range("M100").Value2 = -657.75

# --- Sheet LTW, hunk 47 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value2 = 3405.6155
$ws.Range("I132").Value2 = 3255.4
$ws.Range("J132").Value2 = 3499.5
$ws.Range("K132").Value2 = 9766.2
$ws.Range("L132").Value2 = 10498.5
$ws.Range("M132").Value2 = -7236.200000000001
$ws.Range("N132").Value2 = -15558.5

# --- Sheet WVR, hunk 48 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 379.11765
$ws.Range("I107").Value2 = 302.27274
$ws.Range("K107").Value2 = 906.81822
$ws.Range("M107").Value2 = 1013.18178

